$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" '29.440.52'
Set-TextValue "E2" '  -0.36%  '

Set-TextValue "D3" '1.850.93'
Set-TextValue "E3" '  -0.07%  '

Set-TextValue "D4" '0.9988'
Set-TextValue "E4" '  -0.01%  '

Set-TextValue "D5" '241.03'
Set-TextValue "E5" '  -1.00%  '

Set-TextValue "D6" '0.6334'
Set-TextValue "E6" '  -0.36%  '

Set-TextValue "D7" '0.9998'
Set-TextValue "E7" '  +0.02%  '

Set-TextValue "D8" '3.892.28'
Set-TextValue "E8" '  +108.12%  '

Set-TextValue "B9" 'Dogecoin'
Set-TextValue "C9" 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue "D9" '0.07571'
Set-TextValue "E9" '  +1.25%  '

Set-TextValue "B10" 'WrappedliquidstakedEther2.0'
Set-TextValue "C10" 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue "D10" '4.118.04'
Set-TextValue "E10" '  +91.55%  '

Set-TextValue "D11" '0.2971'
Set-TextValue "E11" '  -1.06%  '

Set-TextValue "B12" 'Solana'
Set-TextValue "C12" 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue "D12" '24.69'
Set-TextValue "E12" '  +1.61%  '

Set-TextValue "B13" 'TRON'
Set-TextValue "C13" 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D13" '0.07725'
Set-TextValue "E13" '  +1.29%  '

Set-TextValue "B14" 'Polkadot'
Set-TextValue "C14" 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue "D14" '4.999'
Set-TextValue "E14" '  -0.67%  '

Set-TextValue "B15" 'Polygon'
Set-TextValue "C15" 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D15" '0.6865'
Set-TextValue "E15" '  +0.05%  '

Set-TextValue "B16" 'Litecoin'
Set-TextValue "C16" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D16" '83.08'
Set-TextValue "E16" '  -0.62%  '

Set-TextValue "B17" 'ShibaInu'
Set-TextValue "C17" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D17" '0.000009955'
Set-TextValue "E17" '  +4.28%  '

Set-TextValue "B18" 'Uniswap'
Set-TextValue "C18" 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D18" '6.214'
Set-TextValue "E18" '  +0.72%  '

Set-TextValue "B19" 'WrappedBTC'
Set-TextValue "C19" 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D19" '29.471.90'
Set-TextValue "E19" '  -0.28%  '

Set-TextValue "B20" 'BitcoinCash'
Set-TextValue "C20" 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D20" '232.30'
Set-TextValue "E20" '  -1.63%  '

Set-TextValue "B21" 'Avalanche'
Set-TextValue "C21" 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D21" '12.50'
Set-TextValue "E21" '  -0.60%  '

Set-TextValue "B22" 'Dai'
Set-TextValue "C22" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D22" '1.0000'
Set-TextValue "E22" '  -0.01%  '

Set-TextValue "B23" 'Chainlink'
Set-TextValue "C23" 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D23" '7.612'
Set-TextValue "E23" '  -1.45%  '

Set-TextValue "B24" 'BinanceUSD'
Set-TextValue "C24" 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue "D24" '0.9999'
Set-TextValue "E24" '  -0.02%  '

Set-TextValue "B25" 'Monero'
Set-TextValue "C25" 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D25" '155.77'
Set-TextValue "E25" '  -0.99%  '

Set-TextValue "B26" 'Stellar'
Set-TextValue "C26" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D26" '0.1390'
Set-TextValue "E26" '  -0.97%  '

Set-TextValue "B27" 'Cosmos'
Set-TextValue "C27" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D27" '8.420'
Set-TextValue "E27" '  -1.00%  '

Set-TextValue "B28" 'EthereumClassic'
Set-TextValue "C28" 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D28" '17.70'
Set-TextValue "E28" '  -0.37%  '

Set-TextValue "B29" 'RocketPoolETH'
Set-TextValue "C29" 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue "D29" '4.095.47'
Set-TextValue "E29" '  +100.09%  '

Set-TextValue "B30" 'PancakeSwap'
Set-TextValue "C30" 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue "D30" '1.469'
Set-TextValue "E30" '  -1.44%  '

Set-TextValue "B31" 'Hedera'
Set-TextValue "C31" 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D31" '0.05831'
Set-TextValue "E31" '  -2.91%  '

Set-TextValue "B32" 'Toncoin'
Set-TextValue "C32" 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue "D32" '1.259'
Set-TextValue "E32" '  +0.72%  '

Set-TextValue "B33" 'Filecoin'
Set-TextValue "C33" 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D33" '4.139'
Set-TextValue "E33" '  +0.42%  '

Set-TextValue "B34" 'InternetComputer(DFINITY)'
Set-TextValue "C34" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D34" '4.026'
Set-TextValue "E34" '  -1.20%  '

Set-TextValue "B35" 'LidoDAOToken'
Set-TextValue "C35" 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue "D35" '1.861'
Set-TextValue "E35" '  -0.57%  '

Set-TextValue "B36" 'ARBITRUM'
Set-TextValue "C36" 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D36" '1.159'
Set-TextValue "E36" '  -1.48%  '

Set-TextValue "B37" 'ImmutableX'
Set-TextValue "C37" 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D37" '0.7177'
Set-TextValue "E37" '  -0.37%  '

Set-TextValue "B38" 'HuobiToken'
Set-TextValue "C38" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D38" '2.595'
Set-TextValue "E38" '  -0.33%  '

Set-TextValue "B39" 'Maker'
Set-TextValue "C39" 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D39" '1.253.53'
Set-TextValue "E39" '  +4.12%  '

Set-TextValue "B40" 'MXToken'
Set-TextValue "C40" 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D40" '2.798'
Set-TextValue "E40" '  -0.13%  '

Set-TextValue "B41" 'VeChain'
Set-TextValue "C41" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D41" '0.01808'
Set-TextValue "E41" '  +1.74%  '

Set-TextValue "B42" 'TrustWalletToken'
Set-TextValue "C42" 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D42" '0.9031'
Set-TextValue "E42" '  -0.72%  '

Set-TextValue "B43" 'FraxShare'
Set-TextValue "C43" 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" '6.114'
Set-TextValue "E43" '  -0.85%  '

Set-TextValue "B44" 'PaxDollar'
Set-TextValue "C44" 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D44" '0.9997'
Set-TextValue "E44" '  +0.03%  '

Set-TextValue "B45" 'Quant'
Set-TextValue "C45" 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D45" '101.74'
Set-TextValue "E45" '  -0.25%  '

Set-TextValue "B46" 'Aave'
Set-TextValue "C46" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D46" '67.16'
Set-TextValue "E46" '  +0.60%  '

Set-TextValue "B47" 'Aptos'
Set-TextValue "C47" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D47" '7.228'
Set-TextValue "E47" '  -0.93%  '

Set-TextValue "B48" 'EnergySwap'
Set-TextValue "C48" 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue "D48" '9.174'
Set-TextValue "E48" '  +0.67%  '

Set-TextValue "B49" 'TheSandbox'
Set-TextValue "C49" 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D49" '0.4024'
Set-TextValue "E49" '  -0.26%  '

Set-TextValue "B50" 'RenderToken'
Set-TextValue "C50" 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D50" '1.694'
Set-TextValue "E50" '  +2.13%  '

Set-TextValue "B51" 'Algorand'
Set-TextValue "C51" 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D51" '0.1128'
Set-TextValue "E51" '  +0.13%  '

